$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 851.3333
$ws.Range("I19").Value = 882.5
$ws.Range("J19").Value = 789
$ws.Range("K19").Value = 882.5
$ws.Range("L19").Value = 789
$ws.Range("M19").Value = -707.5
$ws.Range("N19").Value = -1139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9190.18
$ws.Range("I32").Value = 6388.067
$ws.Range("K32").Value = 6388.067
$ws.Range("M32").Value = -6101.067
$ws.Range("H43").Value = 100000
$ws.Range("J43").Value = 100000
$ws.Range("L43").Value = 100000
$ws.Range("N43").Value = -100626
$ws.Range("H61").Value = 2204.4666
$ws.Range("I61").Value = 2128.3845
$ws.Range("K61").Value = 2128.3845
$ws.Range("M61").Value = -1916.3845
$ws.Range("H132").Value = 1697.7715
$ws.Range("I132").Value = 1609.3448
$ws.Range("K132").Value = 4828.0344
$ws.Range("M132").Value = -2298.0344
$ws.Range("H136").Value = 2204.4666
$ws.Range("I136").Value = 2128.3845
$ws.Range("K136").Value = 6385.1535
$ws.Range("M136").Value = -3835.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3211.2307
$ws.Range("I134").Value = 3261.4285
$ws.Range("J134").Value = 3152.6667
$ws.Range("K134").Value = 9784.2855
$ws.Range("L134").Value = 9458.000100000001
$ws.Range("M134").Value = -7249.2855
$ws.Range("N134").Value = -14528.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 524
$ws.Range("J22").Value = 449
$ws.Range("L22").Value = 449
$ws.Range("N22").Value = -1149
$ws.Range("H31").Value = 5499
$ws.Range("I31").Value = 4479.615
$ws.Range("K31").Value = 4479.615
$ws.Range("M31").Value = -4184.615
$ws.Range("H34").Value = 5499
$ws.Range("I34").Value = 4479.615
$ws.Range("K34").Value = 4479.615
$ws.Range("M34").Value = -4277.615
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21316
$ws.Range("H58").Value = 2932.6072
$ws.Range("I58").Value = 1209.3334
$ws.Range("J58").Value = 4225.0625
$ws.Range("K58").Value = 1209.3334
$ws.Range("L58").Value = 4225.0625
$ws.Range("M58").Value = -1006.3334
$ws.Range("N58").Value = -4631.0625
$ws.Range("H99").Value = 12873.381
$ws.Range("I99").Value = 8032.8
$ws.Range("J99").Value = 17273.908
$ws.Range("K99").Value = 8032.8
$ws.Range("L99").Value = 17273.908
$ws.Range("M99").Value = -6534.8
$ws.Range("N99").Value = -20269.908
$ws.Range("H126").Value = 12873.381
$ws.Range("I126").Value = 8032.8
$ws.Range("J126").Value = 17273.908
$ws.Range("K126").Value = 24098.4
$ws.Range("L126").Value = 51821.724
$ws.Range("M126").Value = -21628.4
$ws.Range("N126").Value = -56761.724
$ws.Range("H136").Value = 2932.6072
$ws.Range("I136").Value = 1209.3334
$ws.Range("J136").Value = 4225.0625
$ws.Range("K136").Value = 3628.0002
$ws.Range("L136").Value = 12675.1875
$ws.Range("M136").Value = -1078.0002
$ws.Range("N136").Value = -17775.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1389.8
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 6000
$ws.Range("N62").Value = -7372
$ws.Range("H65").Value = 1389.8
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 18000
$ws.Range("N65").Value = -24864

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 263.25
$ws.Range("I2").Value = 54.333332
$ws.Range("J2").Value = 388.6
$ws.Range("K2").Value = 54.333332
$ws.Range("L2").Value = 388.6
$ws.Range("M2").Value = 58.666668
$ws.Range("N2").Value = -614.6
$ws.Range("H136").Value = 28793.834
$ws.Range("J136").Value = 28793.834
$ws.Range("L136").Value = 86381.50199999999
$ws.Range("N136").Value = -91481.50199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 11199.4
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 11199.4
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 11199.4
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -11789.4
$ws.Range("H27").Value = 11199.4
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 11199.4
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 11199.4
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -11413.4
$ws.Range("H46").Value = 3490.7273
$ws.Range("J46").Value = 3899.6667
$ws.Range("L46").Value = 3899.6667
$ws.Range("N46").Value = -4275.6667
$ws.Range("H55").Value = 637.0769
$ws.Range("I55").Value = 633.3
$ws.Range("K55").Value = 633.3
$ws.Range("M55").Value = -460.3
$ws.Range("H68").Value = 2279.9167
$ws.Range("I68").Value = 2564.8333
$ws.Range("J68").Value = 1995
$ws.Range("K68").Value = 2564.8333
$ws.Range("L68").Value = 1995
$ws.Range("M68").Value = -1815.8333
$ws.Range("N68").Value = -3493
$ws.Range("H71").Value = 2279.9167
$ws.Range("I71").Value = 2564.8333
$ws.Range("J71").Value = 1995
$ws.Range("K71").Value = 12824.1665
$ws.Range("L71").Value = 9975
$ws.Range("M71").Value = -9080.166499999999
$ws.Range("N71").Value = -17463
$ws.Range("H100").Value = 1433.1666
$ws.Range("I100").Value = 899.75
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 899.75
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -358.75
$ws.Range("N100").Value = -3582
$ws.Range("H124").Value = 57500
$ws.Range("J124").Value = 57500
$ws.Range("L124").Value = 57500
$ws.Range("N124").Value = -67320
$ws.Range("H132").Value = 3553.1316
$ws.Range("I132").Value = 2696.6956
$ws.Range("K132").Value = 8090.0868
$ws.Range("M132").Value = -5560.0868
$ws.Range("H136").Value = 3093.6956
$ws.Range("I136").Value = 3145.1052
$ws.Range("J136").Value = 2849.5
$ws.Range("K136").Value = 9435.3156
$ws.Range("L136").Value = 8548.5
$ws.Range("M136").Value = -6885.3156
$ws.Range("N136").Value = -13648.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1381.125
$ws.Range("I132").Value = 1381.125
$ws.Range("K132").Value = 4143.375
$ws.Range("M132").Value = -1613.375
$ws.Range("H136").Value = 1596.5186
$ws.Range("I136").Value = 1628.28
$ws.Range("J136").Value = 1199.5
$ws.Range("K136").Value = 4884.84
$ws.Range("L136").Value = 1199.5
$ws.Range("M136").Value = -2334.84
$ws.Range("N136").Value = -8698.5
